$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('I3').Value = 'ba'
$ws.Range('J3').Value = 'Appreciation'
$ws.Range('I4').Value = 'sd'
$ws.Range('J4').Value = 'Statement-non-opinion'
$ws.Range('I6').Value = 'aa'
$ws.Range('J6').Value = 'Agree/Accept'
$ws.Range('I24').Value = 'sd'
$ws.Range('J24').Value = 'Statement-non-opinion'
$ws.Range('I38').Value = 'aa'
$ws.Range('J38').Value = 'Agree/Accept'
$ws.Range('I42').Value = 'sd'
$ws.Range('J42').Value = 'Statement-non-opinion'
$ws.Range('I53').Value = '%'
$ws.Range('J53').Value = 'Uninterpretable'
$ws.Range('I55').Value = 'aa'
$ws.Range('J55').Value = 'Agree/Accept'
$ws.Range('I56').Value = 'ba'
$ws.Range('J56').Value = 'Appreciation'
$ws.Range('I60').Value = 'aa'
$ws.Range('J60').Value = 'Agree/Accept'
$ws.Range('I62').Value = 'sd'
$ws.Range('J62').Value = 'Statement-non-opinion'
$ws.Range('I68').Value = 'sd'
$ws.Range('J68').Value = 'Statement-non-opinion'
$ws.Range('I81').Value = 'sd'
$ws.Range('J81').Value = 'Statement-non-opinion'
$ws.Range('I90').Value = 'aa'
$ws.Range('J90').Value = 'Agree/Accept'
$ws.Range('I92').Value = '%'
$ws.Range('J92').Value = 'Uninterpretable'
$ws.Range('I94').Value = 'sd'
$ws.Range('J94').Value = 'Statement-non-opinion'
$ws.Range('I109').Value = '%'
$ws.Range('J109').Value = 'Uninterpretable'
$ws.Range('I117').Value = 'aa'
$ws.Range('J117').Value = 'Agree/Accept'
$ws.Range('I125').Value = 'sd'
$ws.Range('J125').Value = 'Statement-non-opinion'
$ws.Range('I161').Value = 'aa'
$ws.Range('J161').Value = 'Agree/Accept'
$ws.Range('I167').Value = 'sd'
$ws.Range('J167').Value = 'Statement-non-opinion'
$ws.Range('I169').Value = 'sd'
$ws.Range('J169').Value = 'Statement-non-opinion'
$ws.Range('I176').Value = 'sd'
$ws.Range('J176').Value = 'Statement-non-opinion'
$ws.Range('I183').Value = 'sv'
$ws.Range('J183').Value = 'Statement-opinion'
$ws.Range('I195').Value = 'sd'
$ws.Range('J195').Value = 'Statement-non-opinion'
$ws.Range('I196').Value = 'sv'
$ws.Range('J196').Value = 'Statement-opinion'
$ws.Range('I221').Value = 'qy'
$ws.Range('J221').Value = 'Yes-No-Question'
$ws.Range('I224').Value = '%'
$ws.Range('J224').Value = 'Uninterpretable'
$ws.Range('I240').Value = 'sd'
$ws.Range('J240').Value = 'Statement-non-opinion'
$ws.Range('I258').Value = 'aa'
$ws.Range('J258').Value = 'Agree/Accept'
$ws.Range('I271').Value = '%'
$ws.Range('J271').Value = 'Uninterpretable'
$ws.Range('I276').Value = 'sd'
$ws.Range('J276').Value = 'Statement-non-opinion'
$ws.Range('I290').Value = 'sd'
$ws.Range('J290').Value = 'Statement-non-opinion'
$ws.Range('I321').Value = 'sd'
$ws.Range('J321').Value = 'Statement-non-opinion'
$ws.Range('I324').Value = 'sv'
$ws.Range('J324').Value = 'Statement-opinion'
$ws.Range('I331').Value = 'ba'
$ws.Range('J331').Value = 'Appreciation'
$ws.Range('I332').Value = 'ba'
$ws.Range('J332').Value = 'Appreciation'
$ws.Range('I333').Value = 'sv'
$ws.Range('J333').Value = 'Statement-opinion'
$ws.Range('I344').Value = 'sv'
$ws.Range('J344').Value = 'Statement-opinion'
$ws.Range('I356').Value = 'ba'
$ws.Range('J356').Value = 'Appreciation'
$ws.Range('I373').Value = 'sd'
$ws.Range('J373').Value = 'Statement-non-opinion'
$ws.Range('I433').Value = 'sd'
$ws.Range('J433').Value = 'Statement-non-opinion'
$ws.Range('I458').Value = 'aa'
$ws.Range('J458').Value = 'Agree/Accept'
$ws.Range('I462').Value = 'b'
$ws.Range('J462').Value = 'Acknowledge (Backchannel)'
$ws.Range('I463').Value = 'sd'
$ws.Range('J463').Value = 'Statement-non-opinion'
$ws.Range('I467').Value = '%'
$ws.Range('J467').Value = 'Uninterpretable'
$ws.Range('I477').Value = 'sd'
$ws.Range('J477').Value = 'Statement-non-opinion'
$ws.Range('I479').Value = '%'
$ws.Range('J479').Value = 'Uninterpretable'
$ws.Range('I491').Value = 'sd'
$ws.Range('J491').Value = 'Statement-non-opinion'
$ws.Range('I496').Value = 'sd'
$ws.Range('J496').Value = 'Statement-non-opinion'
$ws.Range('I498').Value = 'qy'
$ws.Range('J498').Value = 'Yes-No-Question'
$ws.Range('I504').Value = 'sd'
$ws.Range('J504').Value = 'Statement-non-opinion'
$ws.Range('I507').Value = 'sv'
$ws.Range('J507').Value = 'Statement-opinion'
$ws.Range('I521').Value = 'b'
$ws.Range('J521').Value = 'Acknowledge (Backchannel)'
$ws.Range('I523').Value = 'sd'
$ws.Range('J523').Value = 'Statement-non-opinion'
$ws.Range('I524').Value = 'sd'
$ws.Range('J524').Value = 'Statement-non-opinion'
$ws.Range('I527').Value = 'sv'
$ws.Range('J527').Value = 'Statement-opinion'
$ws.Range('I536').Value = 'sd'
$ws.Range('J536').Value = 'Statement-non-opinion'
$ws.Range('I539').Value = 'sv'
$ws.Range('J539').Value = 'Statement-opinion'
$ws.Range('I544').Value = 'ba'
$ws.Range('J544').Value = 'Appreciation'
$ws.Range('I545').Value = 'sd'
$ws.Range('J545').Value = 'Statement-non-opinion'
$ws.Range('I549').Value = 'sv'
$ws.Range('J549').Value = 'Statement-opinion'
$ws.Range('I560').Value = 'sd'
$ws.Range('J560').Value = 'Statement-non-opinion'
$ws.Range('I561').Value = 'sd'
$ws.Range('J561').Value = 'Statement-non-opinion'
$ws.Range('I562').Value = 'sd'
$ws.Range('J562').Value = 'Statement-non-opinion'
$ws.Range('I568').Value = 'sd'
$ws.Range('J568').Value = 'Statement-non-opinion'
$ws.Range('I577').Value = 'sd'
$ws.Range('J577').Value = 'Statement-non-opinion'
$ws.Range('I578').Value = 'sd'
$ws.Range('J578').Value = 'Statement-non-opinion'
$ws.Range('I585').Value = 'sv'
$ws.Range('J585').Value = 'Statement-opinion'
$ws.Range('I588').Value = 'sv'
$ws.Range('J588').Value = 'Statement-opinion'
$ws.Range('I592').Value = 'sd'
$ws.Range('J592').Value = 'Statement-non-opinion'
$ws.Range('I606').Value = '%'
$ws.Range('J606').Value = 'Uninterpretable'
$ws.Range('I615').Value = 'sd'
$ws.Range('J615').Value = 'Statement-non-opinion'
$ws.Range('I616').Value = 'aa'
$ws.Range('J616').Value = 'Agree/Accept'
$ws.Range('I618').Value = 'aa'
$ws.Range('J618').Value = 'Agree/Accept'
$ws.Range('I620').Value = 'sd'
$ws.Range('J620').Value = 'Statement-non-opinion'
$ws.Range('I623').Value = 'sd'
$ws.Range('J623').Value = 'Statement-non-opinion'
$ws.Range('I634').Value = 'aa'
$ws.Range('J634').Value = 'Agree/Accept'
$ws.Range('I641').Value = 'sd'
$ws.Range('J641').Value = 'Statement-non-opinion'
$ws.Range('I650').Value = 'aa'
$ws.Range('J650').Value = 'Agree/Accept'
$ws.Range('I654').Value = 'aa'
$ws.Range('J654').Value = 'Agree/Accept'
